$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Charts")

# Add the "coming soon" message to the Charts tab and make it the active sheet/cell.
$ws.Range("A1").Value = "Automatically generated chart(s) coming soon to this tab."
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
